$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness values (column C, rows 2-252) from 7293 to 7590
$ws.Range("C2:C252").Value = 7590
